$d = $word.ActiveDocument

# Find the paragraph with the authors' names
$findRng = $d.Content
$found = $findRng.Find.Execute("Edison Achalma, Yeno Areste, y Cristían Galindo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Host "ERROR: target paragraph not found"
} else {
    # Determine the 1-based index of the found paragraph
    $targetIndex = $d.Range(0, $findRng.Start).Paragraphs.Count + 1

    # Insert a paragraph break right after the found text
    $insertPos = $findRng.End
    $d.Range($insertPos, $insertPos).InsertAfter("`r")

    # The newly created paragraph is the one right after the target paragraph
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Style = "Author"
    $newPara.Range.InsertBefore("Economía, Universidad Nacional de San Cristóbal de Huamanga")

    Write-Host "Inserted new Author paragraph after paragraph index $targetIndex"
}
